# Upload Leave Card 12/27/2023 4:01 PM
# Adds a VL leave usage entry (1 day, used 11/23/2023) and a second VL
# leave usage entry (3 days, used 12/27-29/2023), together with the SL
# and SP entries earned in between, to the Sheet1 (Leave Card) table,
# and inserts a "2024" year-divider row ahead of the 2024 periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# --- 1. Insert two new blank rows inside the table body -----------------
# Row 210 is inserted first (pushes old row 210 -> 211, old 211 -> 212, ...).
$ws.Rows("210:210").Insert()
# Give the freshly inserted row the same look (borders / number formats /
# styles) as the row above it - Insert() on its own drops formatting.
$ws.Range("A209:K209").Copy()
$ws.Range("A210:K210").PasteSpecial(-4122)

# Row 212 is inserted next (pushes old row 211 [now at row 211] -> 213, ...).
$ws.Rows("212:212").Insert()
$ws.Range("A211:K211").Copy()
$ws.Range("A212:K212").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Resize the table to include the two new rows ----------------------
$tbl.Resize($ws.Range("A8:K344"))

# --- 3. Re-create the "EARNED " helper-column formula on the new rows -----
# (Rows.Insert() does not propagate the table's calculated column formula.)
$ws.Range("G210").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G212").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 4. Fill in the leave entries ------------------------------------------

# Row 208 : period 10/01/2023 - SL(1-0-0) earned 1.25, used 10/16/2023
$ws.Range("B208").Value = "SL(1-0-0)"
$ws.Range("C207").Copy()
$ws.Range("C208").PasteSpecial(-4122)
$ws.Range("C208").Value = 1.25
$ws.Range("H208").Value = 1
$ws.Range("K205").Copy()
$ws.Range("K208").PasteSpecial(-4122)
$ws.Range("K208").Value = 45215

# Row 209 : period 11/01/2023 - SP(1-0-0) earned 1.25, used 12/01/2023
$ws.Range("B209").Value = "SP(1-0-0)"
$ws.Range("C206").Copy()
$ws.Range("C209").PasteSpecial(-4122)
$ws.Range("C209").Value = 1.25
$ws.Range("K206").Copy()
$ws.Range("K209").PasteSpecial(-4122)
$ws.Range("K209").Value = 45261

# Row 210 (new) : VL(1-0-0), 1 day absence w/ pay, used 11/23/2023
$ws.Range("A210").Value = ""
$ws.Range("B210").Value = "VL(1-0-0)"
$ws.Range("D210").Value = 1
$ws.Range("K207").Copy()
$ws.Range("K210").PasteSpecial(-4122)
$ws.Range("K210").Value = 45253

# Row 211 (old row 210, shifted down) : period 12/01/2023, VL(3-0-0),
# 3 days absence w/ pay, remarks = "12/27-29/2023"
$ws.Range("B211").Value = "VL(3-0-0)"
$ws.Range("D211").Value = 3

# Row 212 (new) : "2024" year divider - set up *before* the K211 remark
# below so the two brand-new shared strings land in the same order as
# the source edit ("2024" then "12/27-29/2023").
$ws.Range("A212").NumberFormat = "@"
$ws.Range("A212").Value = "2024"
$ws.Range("A196").Copy()
$ws.Range("A212").PasteSpecial(-4122)

$ws.Range("K211").Value = "12/27-29/2023"

# --- 5. Cosmetic: move the saved selection like the author's session -----
$ws.Range("K211").Select()

$wb.Save()
